$d = $word.ActiveDocument

$BLUE  = 12611584   # wdColor value equivalent to RGB hex 0070C0
$GREEN = 5287936    # wdColor value equivalent to RGB hex 00B050

# ---------------------------------------------------------------------
# helper: strip any leftover w:rsid* baggage from a Range by
# round-tripping its text through a throwaway placeholder, then apply
# the requested font color. (Just setting Font.Color on a Range leaves
# any pre-existing w:rsid* attribute in place; actually replacing the
# text forces a fresh run with no rsid - which is how the runs in the
# target revision were re-serialized.)
# ---------------------------------------------------------------------
function Recolor-Range {
    param($rng, $color)
    $orig = $rng.Text
    $rng.Text = "^"
    $rng.Text = $orig
    $rng.Font.Color = $color
}

function Append-ColoredText {
    param($para, $text, $color)
    $r = $para.Range
    $insertStart = $r.End - 1
    $ins = $d.Range($insertStart, $insertStart)
    $ins.InsertAfter($text)
    # the freshly-inserted Range otherwise clones the rsid of whatever
    # run precedes it; round-trip the text so it serializes as a clean
    # run with no w:rsid* attribute, matching the target revision.
    Recolor-Range $ins $color
}

# ---------------------------------------------------------------------
# 1) "Add "new user registration" link to homepage header " -> append
#    "   (DAN)" in blue (0070C0)
# ---------------------------------------------------------------------
$p1 = $d.Paragraphs.Item(183)
Append-ColoredText $p1 "   (DAN)" $BLUE

# ---------------------------------------------------------------------
# 2) "...and most all other pages" -> append "   (DAN)" in blue
# ---------------------------------------------------------------------
$p2 = $d.Paragraphs.Item(184)
Append-ColoredText $p2 "   (DAN)" $BLUE

# ---------------------------------------------------------------------
# 3) DOCKER paragraph -> whole paragraph (incl. end-of-paragraph mark)
#    turns green, then a *separate* "    (GREG)" run is appended, also
#    green.
# ---------------------------------------------------------------------
$p3 = $d.Paragraphs.Item(215)
$p3.Range.Font.Color = $GREEN
$r3 = $p3.Range
$insertStart3 = $r3.End - 1
$ins3 = $d.Range($insertStart3, $insertStart3)
$ins3.InsertAfter("    (GREG)")
# Toggle an unrelated property so the host keeps this as its own run
# instead of silently merging it into the identically-colored "DOCKER"
# run that precedes it (mirrors the two separate <w:r> in the target).
$ins3.Bold = 1
$ins3.Font.Color = $GREEN
$ins3.Bold = 0

# ---------------------------------------------------------------------
# 4) EnvironmentVariables paragraph -> whole paragraph (incl. mark)
#    turns green; the four existing runs (split around the proofErr
#    spell-check markers) are each individually recolored and their
#    rsid stripped, but stay separate runs.
# ---------------------------------------------------------------------
$p4 = $d.Paragraphs.Item(216)
$base4 = $p4.Range.Start

$seg4a = $d.Range($base4, $base4 + 20)         # "EnvironmentVariables"
Recolor-Range $seg4a $GREEN

$seg4b = $d.Range($base4 + 20, $base4 + 26)    # " setup"
Recolor-Range $seg4b $GREEN

$seg4c = $d.Range($base4 + 26, $base4 + 32)    # " into "
Recolor-Range $seg4c $GREEN

$seg4d = $d.Range($base4 + 32, $base4 + 38)    # "docker"
Recolor-Range $seg4d $GREEN

$p4.Range.Font.Color = $GREEN

# ---------------------------------------------------------------------
# 5) Project1_UName paragraph -> whole paragraph green; "Project1_UName"
#    and "    - " merge into a single run (both already share identical
#    non-rsid formatting), "postgres" (wrapped in proofErr) stays its
#    own run.
# ---------------------------------------------------------------------
$p5 = $d.Paragraphs.Item(217)
$base5 = $p5.Range.Start

$seg5a = $d.Range($base5, $base5 + 20)         # "Project1_UName    - "
Recolor-Range $seg5a $GREEN

$seg5b = $d.Range($base5 + 20, $p5.Range.End - 1)   # "postgres"
Recolor-Range $seg5b $GREEN

$p5.Range.Font.Color = $GREEN

# ---------------------------------------------------------------------
# 6) Project1_UPassword paragraph -> whole paragraph green; all 3 runs
#    merge into a single run.
# ---------------------------------------------------------------------
$p6 = $d.Paragraphs.Item(218)
$r6 = $p6.Range
$r6.MoveEnd(1, -1)
Recolor-Range $r6 $GREEN
$p6.Range.Font.Color = $GREEN

# ---------------------------------------------------------------------
# 7) Project1_URL paragraph -> whole paragraph green; all 3 runs merge
#    into a single run.
# ---------------------------------------------------------------------
$p7 = $d.Paragraphs.Item(219)
$r7 = $p7.Range
$r7.MoveEnd(1, -1)
Recolor-Range $r7 $GREEN
$p7.Range.Font.Color = $GREEN
